$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '42.096.98'
$ws.Cells.Item(2, 5).Value = '  +5.63%  '

$ws.Cells.Item(3, 4).Value = '2.267.02'

$ws.Cells.Item(4, 4).Value = '0.999'
$ws.Cells.Item(4, 5).Value = '  -0.08%  '

$ws.Cells.Item(5, 4).Value = '302.69'
$ws.Cells.Item(5, 5).Value = '  +3.80%  '

$ws.Cells.Item(6, 4).Value = '93.24'
$ws.Cells.Item(6, 5).Value = '  +7.34%  '

$ws.Cells.Item(7, 4).Value = '0.535'
$ws.Cells.Item(7, 5).Value = '  +4.80%  '

$ws.Cells.Item(8, 5).Value = '  -0.05%  '

$ws.Cells.Item(9, 4).Value = '0.487'
$ws.Cells.Item(9, 5).Value = '  +4.20%  '

$ws.Cells.Item(10, 4).Value = '32.93'
$ws.Cells.Item(10, 5).Value = '  +7.96%  '

$ws.Cells.Item(11, 4).Value = '54.82'
$ws.Cells.Item(11, 5).Value = '  +9.67%  '

$ws.Cells.Item(12, 5).Value = '  +3.00%  '

$ws.Cells.Item(13, 5).Value = '  +3.41%  '

$ws.Cells.Item(14, 4).Value = '6.72'
$ws.Cells.Item(14, 5).Value = '  +4.47%  '

$ws.Cells.Item(15, 4).Value = '2.619.45'
$ws.Cells.Item(15, 5).Value = '  +2.33%  '

$ws.Cells.Item(16, 4).Value = '14.20'
$ws.Cells.Item(16, 5).Value = '  +3.64%  '

$ws.Cells.Item(17, 4).Value = '2.267.70'
$ws.Cells.Item(17, 5).Value = '  +0.31%  '

$ws.Cells.Item(18, 4).Value = '0.759'
$ws.Cells.Item(18, 5).Value = '  +4.03%  '

$ws.Cells.Item(19, 4).Value = '41.972.86'
$ws.Cells.Item(19, 5).Value = '  +5.46%  '

$ws.Cells.Item(20, 5).Value = '  +9.07%  '

$ws.Cells.Item(21, 4).Value = '0.0₃0910'
$ws.Cells.Item(21, 5).Value = '  +3.01%  '

$ws.Cells.Item(22, 4).Value = '5.97'
$ws.Cells.Item(22, 5).Value = '  +3.80%  '

$ws.Cells.Item(23, 4).Value = '67.47'
$ws.Cells.Item(23, 5).Value = '  +2.96%  '

$ws.Cells.Item(24, 4).Value = '242.89'
$ws.Cells.Item(24, 5).Value = '  +2.67%  '

$ws.Cells.Item(25, 5).Value = '  +6.09%  '

$ws.Cells.Item(26, 4).Value = '1.00'
$ws.Cells.Item(26, 5).Value = '  +0.05%  '

$ws.Cells.Item(27, 5).Value = '  +4.92%  '

$ws.Cells.Item(28, 4).Value = '23.99'

$ws.Cells.Item(29, 4).Value = '2.18'
$ws.Cells.Item(29, 5).Value = '  +1.77%  '

$ws.Cells.Item(30, 5).Value = '  +5.49%  '

$ws.Cells.Item(31, 4).Value = '34.19'
$ws.Cells.Item(31, 5).Value = '  +7.78%  '

$ws.Cells.Item(32, 4).Value = '158.46'
$ws.Cells.Item(32, 5).Value = '  +1.07%  '

$ws.Cells.Item(33, 5).Value = '  -0.07%  '

$ws.Cells.Item(34, 5).Value = '  +5.02%  '

$ws.Cells.Item(35, 2).Value = 'Hedera'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(35, 4).Value = '0.0745'
$ws.Cells.Item(35, 5).Value = '  +5.05%  '

$ws.Cells.Item(36, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Cells.Item(36, 4).Value = '3.08'
$ws.Cells.Item(36, 5).Value = '  +6.04%  '

$ws.Cells.Item(37, 5).Value = '  +3.20%  '

$ws.Cells.Item(38, 5).Value = '  +6.89%  '

$ws.Cells.Item(39, 4).Value = '16.68'
$ws.Cells.Item(39, 5).Value = '  +9.85%  '

$ws.Cells.Item(40, 5).Value = '  +4.26%  '

$ws.Cells.Item(41, 5).Value = '  +5.60%  '

$ws.Cells.Item(42, 4).Value = '3.98'
$ws.Cells.Item(42, 5).Value = '  +7.07%  '

$ws.Cells.Item(43, 4).Value = '20.11'
$ws.Cells.Item(43, 5).Value = '  +13.84%  '

$ws.Cells.Item(44, 4).Value = '2.052.85'
$ws.Cells.Item(44, 5).Value = '  -2.69%  '

$ws.Cells.Item(45, 4).Value = '0.0281'
$ws.Cells.Item(45, 5).Value = '  +4.75%  '

$ws.Cells.Item(46, 4).Value = '10.12'
$ws.Cells.Item(46, 5).Value = '  +2.21%  '

$ws.Cells.Item(47, 4).Value = '2.94'
$ws.Cells.Item(47, 5).Value = '  +8.82%  '

$ws.Cells.Item(48, 5).Value = '  -4.37%  '

$ws.Cells.Item(49, 4).Value = '2.491.66'
$ws.Cells.Item(49, 5).Value = '  +2.35%  '

$ws.Cells.Item(50, 5).Value = '  +2.82%  '

$ws.Cells.Item(51, 2).Value = 'MultiversX'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Cells.Item(51, 4).Value = '52.07'
$ws.Cells.Item(51, 5).Value = '  +6.47%  '
